# Apply the changes described by the commit:
#  1. Fix the typo in the laboratory name "GLICOL Y CYA" -> "GLICOL Y CIA"
#     (shared across every data row in column E).
#  2. Update the saved view/selection: drop the scrolled-down
#     topLeftCell and move the active selection to H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct "GLICOL Y CYA" -> "GLICOL Y CIA" for every row that uses it ---
$ws.Range("E2:E145").Value = "GLICOL Y CIA"

# --- 2. Reset the view so the selection sits on H7 (clears topLeftCell scroll) ---
$ws.Range("H7").Select()
